$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings) ---
$ws.Range("A8").Value = "Volume 33   Number  4"
$ws.Range("C9").Value = "Report Covering the Week  1/19/2026  Through  1/25/2026"

# --- Column H width adjustment (widened to fit new data) ---
$ws.Columns("H").ColumnWidth = 6.71

# --- Simple numeric value updates (style unchanged) ---
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("F16").Value = 2
$ws.Range("H16").Value = 100
$ws.Range("M16").Value = -60
$ws.Range("N16").Value = -95.555555555555
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = 600
$ws.Range("L17").Value = -30
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = -56.25
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 2
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = -33.333333333333
$ws.Range("I18").Value = 5
$ws.Range("J18").Value = 9
$ws.Range("K18").Value = -44.444444444444
$ws.Range("L18").Value = -50
$ws.Range("M18").Value = -37.5
$ws.Range("N18").Value = -94.117647058823
$ws.Range("C19").Value = 11
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 44
$ws.Range("G19").Value = 40
$ws.Range("H19").Value = 10
$ws.Range("I19").Value = 39
$ws.Range("J19").Value = 37
$ws.Range("K19").Value = 5.405405405405
$ws.Range("L19").Value = -18.75
$ws.Range("M19").Value = -15.217391304347
$ws.Range("N19").Value = -68.292682926829
$ws.Range("F20").Value = 3
$ws.Range("H20").Value = 50
$ws.Range("N20").Value = -93.023255813953
$ws.Range("C21").Value = 13
$ws.Range("D21").Value = 13
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 63
$ws.Range("G21").Value = 54
$ws.Range("H21").Value = 16.666666666666
$ws.Range("I21").Value = 56
$ws.Range("J21").Value = 51
$ws.Range("K21").Value = 9.803921568627
$ws.Range("L21").Value = -27.272727272727
$ws.Range("M21").Value = -15.151515151515
$ws.Range("N21").Value = -82.108626198083
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 100
$ws.Range("L22").Value = -33.333333333333
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 32
$ws.Range("E24").Value = -50
$ws.Range("F24").Value = 60
$ws.Range("G24").Value = 95
$ws.Range("H24").Value = -36.842105263157
$ws.Range("I24").Value = 52
$ws.Range("J24").Value = 91
$ws.Range("K24").Value = -42.857142857142
$ws.Range("L24").Value = -8.771929824561
$ws.Range("M24").Value = 23.809523809523
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 26
$ws.Range("E25").Value = -57.692307692307
$ws.Range("F25").Value = 42
$ws.Range("G25").Value = 66
$ws.Range("H25").Value = -36.363636363636
$ws.Range("I25").Value = 39
$ws.Range("J25").Value = 66
$ws.Range("K25").Value = -40.90909090909
$ws.Range("L25").Value = -9.302325581395
$ws.Range("F26").Value = 6
$ws.Range("G26").Value = 12
$ws.Range("H26").Value = -50
$ws.Range("I26").Value = 5
$ws.Range("J26").Value = 7
$ws.Range("K26").Value = -28.571428571428
$ws.Range("L26").Value = -70.588235294117
$ws.Range("M26").Value = -72.222222222222
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 0
$ws.Range("I28").Value = 2
$ws.Range("K28").Value = 100

# --- Cells requiring a type/style change (number <-> text placeholder) ---
# Donor cells with stable styles, never modified by this edit:
#   style 13 (text/General) donor -> C23
#   style 14 (#,##0 integer) donor -> F15
#   style 15 (#,##0.0 percent) donor -> K15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0"
$ws.Range("C23").Copy()
$ws.Range("D15").PasteSpecial(-4122)

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "***.*"
$ws.Range("C23").Copy()
$ws.Range("E15").PasteSpecial(-4122)

$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "0"
$ws.Range("C23").Copy()
$ws.Range("C16").PasteSpecial(-4122)

$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "0"
$ws.Range("C23").Copy()
$ws.Range("C17").PasteSpecial(-4122)

$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "0"
$ws.Range("C23").Copy()
$ws.Range("C20").PasteSpecial(-4122)

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0"
$ws.Range("C23").Copy()
$ws.Range("D20").PasteSpecial(-4122)

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "***.*"
$ws.Range("C23").Copy()
$ws.Range("E20").PasteSpecial(-4122)

$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("C23").Copy()
$ws.Range("C22").PasteSpecial(-4122)

$ws.Range("D22").Value = 1
$ws.Range("F15").Copy()
$ws.Range("D22").PasteSpecial(-4122)

$ws.Range("E22").Value = -100
$ws.Range("K15").Copy()
$ws.Range("E22").PasteSpecial(-4122)

$ws.Range("J22").Value = 1
$ws.Range("F15").Copy()
$ws.Range("J22").PasteSpecial(-4122)

$ws.Range("K22").Value = 100
$ws.Range("K15").Copy()
$ws.Range("K22").PasteSpecial(-4122)

$ws.Range("C26").Value = 1
$ws.Range("F15").Copy()
$ws.Range("C26").PasteSpecial(-4122)

$ws.Range("D26").Value = 3
$ws.Range("F15").Copy()
$ws.Range("D26").PasteSpecial(-4122)

$ws.Range("E26").Value = -66.666666666666
$ws.Range("K15").Copy()
$ws.Range("E26").PasteSpecial(-4122)

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("C23").Copy()
$ws.Range("D27").PasteSpecial(-4122)

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("C23").Copy()
$ws.Range("E27").PasteSpecial(-4122)

$ws.Range("D31").Value = 2
$ws.Range("F15").Copy()
$ws.Range("D31").PasteSpecial(-4122)

$ws.Range("E31").Value = -100
$ws.Range("K15").Copy()
$ws.Range("E31").PasteSpecial(-4122)

$ws.Range("G31").Value = 2
$ws.Range("F15").Copy()
$ws.Range("G31").PasteSpecial(-4122)

$ws.Range("H31").Value = -100
$ws.Range("K15").Copy()
$ws.Range("H31").PasteSpecial(-4122)

$ws.Range("J31").Value = 2
$ws.Range("F15").Copy()
$ws.Range("J31").PasteSpecial(-4122)

$ws.Range("K31").Value = -100
$ws.Range("K15").Copy()
$ws.Range("K31").PasteSpecial(-4122)

$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "0"
$ws.Range("C23").Copy()
$ws.Range("C33").PasteSpecial(-4122)

$excel.CutCopyMode = $false